$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "88-27="
$t.Cell(1, 2).Range.Text = "59-41="
$t.Cell(1, 3).Range.Text = "53-37="
$t.Cell(1, 4).Range.Text = "87-2="
$t.Cell(1, 5).Range.Text = "71-18="

$t.Cell(2, 1).Range.Text = "35+14="
$t.Cell(2, 2).Range.Text = "15+41="
$t.Cell(2, 3).Range.Text = "28-4="
$t.Cell(2, 4).Range.Text = "79-27="
$t.Cell(2, 5).Range.Text = "59-44="

$t.Cell(3, 1).Range.Text = "1+36="
$t.Cell(3, 2).Range.Text = "26+11="
$t.Cell(3, 3).Range.Text = "12+1="
$t.Cell(3, 4).Range.Text = "51-38="
$t.Cell(3, 5).Range.Text = "57+14="

$t.Cell(4, 1).Range.Text = "74+3="
$t.Cell(4, 2).Range.Text = "38+54="
$t.Cell(4, 3).Range.Text = "88-62="
$t.Cell(4, 4).Range.Text = "52-36="
$t.Cell(4, 5).Range.Text = "16+31="

$t.Cell(5, 1).Range.Text = "88-46="
$t.Cell(5, 2).Range.Text = "65-2="
$t.Cell(5, 3).Range.Text = "22+11="
$t.Cell(5, 4).Range.Text = "58-26="
$t.Cell(5, 5).Range.Text = "99-21="

$t.Cell(6, 1).Range.Text = "73-72="
$t.Cell(6, 2).Range.Text = "83+8="
$t.Cell(6, 3).Range.Text = "42-7="
$t.Cell(6, 4).Range.Text = "13+66="
$t.Cell(6, 5).Range.Text = "61-13="

$t.Cell(7, 1).Range.Text = "38+51="
$t.Cell(7, 2).Range.Text = "40+23="
$t.Cell(7, 3).Range.Text = "94-60="
$t.Cell(7, 4).Range.Text = "77-44="
$t.Cell(7, 5).Range.Text = "22-8="

$t.Cell(8, 1).Range.Text = "23-10="
$t.Cell(8, 2).Range.Text = "49+16="
$t.Cell(8, 3).Range.Text = "50+9="
$t.Cell(8, 4).Range.Text = "16+48="
$t.Cell(8, 5).Range.Text = "36+58="

$t.Cell(9, 1).Range.Text = "93-34="
$t.Cell(9, 2).Range.Text = "87-71="
$t.Cell(9, 3).Range.Text = "82-19="
$t.Cell(9, 4).Range.Text = "52+38="
$t.Cell(9, 5).Range.Text = "9+65="

$t.Cell(10, 1).Range.Text = "39+54="
$t.Cell(10, 2).Range.Text = "81+14="
$t.Cell(10, 3).Range.Text = "98-9="
$t.Cell(10, 4).Range.Text = "16-1="
$t.Cell(10, 5).Range.Text = "47+16="

$t.Cell(11, 1).Range.Text = "80-61="
$t.Cell(11, 2).Range.Text = "73-36="
$t.Cell(11, 3).Range.Text = "31-18="
$t.Cell(11, 4).Range.Text = "7+1="
$t.Cell(11, 5).Range.Text = "82-20="

$t.Cell(12, 1).Range.Text = "14+36="
$t.Cell(12, 2).Range.Text = "15+12="
$t.Cell(12, 3).Range.Text = "5+53="
$t.Cell(12, 4).Range.Text = "85+5="
$t.Cell(12, 5).Range.Text = "54+8="

$t.Cell(13, 1).Range.Text = "13+84="
$t.Cell(13, 2).Range.Text = "99-46="
$t.Cell(13, 3).Range.Text = "99-18="
$t.Cell(13, 4).Range.Text = "19+1="
$t.Cell(13, 5).Range.Text = "3+30="

$t.Cell(14, 1).Range.Text = "49+38="
$t.Cell(14, 2).Range.Text = "20+65="
$t.Cell(14, 3).Range.Text = "7+25="
$t.Cell(14, 4).Range.Text = "71-2="
$t.Cell(14, 5).Range.Text = "24-11="

$t.Cell(15, 1).Range.Text = "69-68="
$t.Cell(15, 2).Range.Text = "9+30="
$t.Cell(15, 3).Range.Text = "12+25="
$t.Cell(15, 4).Range.Text = "13+75="
$t.Cell(15, 5).Range.Text = "0+37="

$t.Cell(16, 1).Range.Text = "82-61="
$t.Cell(16, 2).Range.Text = "88-16="
$t.Cell(16, 3).Range.Text = "84-40="
$t.Cell(16, 4).Range.Text = "18-17="
$t.Cell(16, 5).Range.Text = "75-49="

$t.Cell(17, 1).Range.Text = "19+20="
$t.Cell(17, 2).Range.Text = "93+2="
$t.Cell(17, 3).Range.Text = "94-86="
$t.Cell(17, 4).Range.Text = "55-23="
$t.Cell(17, 5).Range.Text = "30+9="

$t.Cell(18, 1).Range.Text = "35+41="
$t.Cell(18, 2).Range.Text = "50+30="
$t.Cell(18, 3).Range.Text = "98+0="
$t.Cell(18, 4).Range.Text = "76+6="
$t.Cell(18, 5).Range.Text = "65+22="

$t.Cell(19, 1).Range.Text = "62-3="
$t.Cell(19, 2).Range.Text = "93-90="
$t.Cell(19, 3).Range.Text = "30+25="
$t.Cell(19, 4).Range.Text = "83-76="
$t.Cell(19, 5).Range.Text = "44-10="

$t.Cell(20, 1).Range.Text = "82+17="
$t.Cell(20, 2).Range.Text = "98-46="
$t.Cell(20, 3).Range.Text = "10+10="
$t.Cell(20, 4).Range.Text = "61+28="
$t.Cell(20, 5).Range.Text = "59+20="
